$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.718.38"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.480.64"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.49"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.66"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.485.22"
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.076.74"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.39"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.491.57"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.829.11"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.90"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.28"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.14"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.624.52"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.09"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  +3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.55"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("E32").Value = "  -1.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.491.79"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.44"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.32"
$ws.Range("E37").Value = "  +5.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "160.11"
$ws.Range("E40").Value = "  -4.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.808"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.20"
$ws.Range("E43").Value = "  +4.42%  "

$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.68"
$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.84"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.415.16"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.900"
$ws.Range("E51").Value = "  +2.11%  "
